$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Add the two new rows of data (shared strings match: values stored as text)
$ws.Range("A80").Value = "FSZZ108"
$ws.Range("B80").Value = "917.19"
$ws.Range("A81").Value = "FSZZ106"
$ws.Range("B81").Value = "1078.19"

# Update selection / view state to match the diff
$ws.Range("B81").Select()
$excel.ActiveWindow.ScrollRow = 61

$excel.ActiveWindow.WindowState = -4143
